# STAT 4110 schedule update — "updated lecture slides for chapter 3"
#
# The Chapter 3 block (weeks of Mar 2 - Apr 3, rows 9-14) is reshuffled:
# the quiz labels are simplified/renumbered, an "Assignment 2 Review"
# session is introduced, and the 8 Chapter-3 lectures are re-spread
# across rows 9-12. Row 15 also loses the old "Final Project
# Presentations" placeholders (now just the Chapter 5 lectures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text values (Chapter 3 lecture reorder + quiz label updates) ---
$ws.Range("B5").Value  = 'Quiz 1'
$ws.Range("B9").Value  = 'Quiz 2'
$ws.Range("D9").Value  = 'Lecture 1 (Chapter 3)'
$ws.Range("E9").Value  = 'Project Summary & Team Plan Due Mar 1, Assignment 2 Published Mar 1'
$ws.Range("B10").Value = 'Assignment 2 Review'
$ws.Range("C10").Value = 'Lecture 2 (Chapter 3)'
$ws.Range("D10").Value = 'Lecture 3 (Chapter 3)'
$ws.Range("B11").Value = 'Lecture 4 (Chapter 3)'
$ws.Range("C11").Value = 'Lecture 5 (Chapter 3)'
$ws.Range("D11").Value = 'Lecture 6 (Chapter 3)'
$ws.Range("E11").Value = 'Assignment 2 Due Mar 20'
$ws.Range("B12").Value = 'Lecture 7 (Chapter 3)'
$ws.Range("C12").Value = 'Lecture 8 (Chapter 3)'
$ws.Range("D12").Value = 'Homework Review Chapter 3'
$ws.Range("B13").Value = 'Quiz 3 '
$ws.Range("C13").Value = 'Lecture 1 (Chapter 4)'
$ws.Range("D13").Value = 'Lecture 2 (Chapter 4)'
$ws.Range("B14").Value = 'Homework Review Chapter 4'
$ws.Range("C14").Value = 'Quiz 4 '
$ws.Range("D14").Value = 'Lecture 1 (Chapter 5)'
$ws.Range("B15").Value = 'Lecture 2 (Chapter 5)'
$ws.Range("C15").Value = 'Lecture 3 (Chapter 5)'

# --- Update bold formatting to match the reshuffled content ---
$ws.Range("B9").Font.Bold  = $true
$ws.Range("B10").Font.Bold = $true
$ws.Range("B12").Font.Bold = $false
$ws.Range("C12").Font.Bold = $false
$ws.Range("D12").Font.Bold = $true
$ws.Range("B13").Font.Bold = $true
$ws.Range("C13").Font.Bold = $false
$ws.Range("D13").Font.Bold = $false
$ws.Range("B14").Font.Bold = $true
$ws.Range("C14").Font.Bold = $true

# --- Row 1 and row 15 no longer need the taller (wrapped) 30pt row height ---
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(15).AutoFit()

# --- Update the saved cursor/selection position ---
$ws.Range("E14").Select()
